# Re-colour the deck's theme (ppt/theme/theme1.xml) from the "Integral /
# Red Violet" scheme to the default "Office" colour scheme (the scheme
# that lives in theme2.xml, used by the Notes Master, in the source
# file). Font scheme and format scheme are already identical between
# the two themes, so only the 12 theme colours need to change.

function ToRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index map (confirmed empirically): 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1
# 6=accent2 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
$tcs.Item(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink 954F72
